$d = $word.ActiveDocument

# 1. Update "Curso (semestre ideal)" line
$d.Content.Find.Execute("Curso (semestre ideal): EQN (12)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Curso (semestre ideal): EQD (10), EQN (12)", 2)

# 2. Remove the "Requisitos" heading paragraph and the requirement bullet paragraph that follows it
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    if ($text -match "Requisitos" -or $text -match "LOQ4044") {
        $p.Range.Delete()
    }
}
